# Apply the 2021-01-19 data update to the "Fonds de solidarite - volet 2"
# sheet. All data cells in this sheet are stored as text (inlineStr), so
# every new value is written with a leading apostrophe to force Excel to
# keep it as text instead of auto-converting it to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 10
$ws.Range("C10").Value = "'500"
$ws.Range("D10").Value = "'442"
$ws.Range("E10").Value = "'2946256.06"

# row 11
$ws.Range("C11").Value = "'233"
$ws.Range("D11").Value = "'199"
$ws.Range("E11").Value = "'2173031.14"

# row 12
$ws.Range("C12").Value = "'73"
$ws.Range("E12").Value = "'1093249.89"

# row 13
$ws.Range("C13").Value = "'24"
$ws.Range("E13").Value = "'578777.00"

# row 30
$ws.Range("C30").Value = "'180"
$ws.Range("E30").Value = "'788748.58"

# row 31
$ws.Range("C31").Value = "'76"
$ws.Range("D31").Value = "'72"
$ws.Range("E31").Value = "'468729.02"

# row 32
$ws.Range("C32").Value = "'17"
$ws.Range("E32").Value = "'125000.00"

# row 74
$ws.Range("C74").Value = "'5"
$ws.Range("E74").Value = "'185000.00"

# row 93
$ws.Range("C93").Value = "'1116"
$ws.Range("D93").Value = "'1012"
$ws.Range("E93").Value = "'6102866.71"

# row 95
$ws.Range("C95").Value = "'194"
$ws.Range("E95").Value = "'2315043.41"
